# Adds two new "CHAMPIONS" stat sheets by duplicating the existing Sheet1
# (hitting) and Sheet2 (pitching) sheets - this carries over all cell
# values/types/styles and the per-row hyperlinks automatically - then
# tweaks the handful of stat values that differ on the new sheets and
# restores the expected tab/selection state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# New sheet order ends up: Sheet1, Sheet2, CHAMPIONS Hitting, CHAMPIONS Pitching
$ws1.Copy($null, $ws2)
$hitting = $wb.Worksheets.Item($ws2.Index + 1)
$hitting.Name = "CHAMPIONS Hitting"

$ws2.Copy($null, $hitting)
$pitching = $wb.Worksheets.Item($hitting.Index + 1)
$pitching.Name = "CHAMPIONS Pitching"

# CHAMPIONS Hitting has a few stats that differ from Sheet1's numbers
$hitting.Range("B5").Value = 641
$hitting.Range("E5").Value = 88
$hitting.Range("F15").Value = 0.269

# CHAMPIONS Pitching has a few stats that differ from Sheet2's numbers
$pitching.Range("B11").Value = 106
$pitching.Range("D15").Value = 3.16

# Restore each sheet's selection / active-cell state
$ws1.Activate()
$ws1.Range("G15").Select()

$ws2.Activate()
$ws2.Range("G15").Select()

$pitching.Activate()
$pitching.Range("B12").Select()

# CHAMPIONS Hitting ends up as the active tab
$hitting.Activate()
$hitting.Range("I12").Select()
